# Fruta / hortaliza, semanal
# Rotates the weekly price-report rows (2-9) one slot: the data that used
# to live in row 9 now lives in row 2, row 8 -> row 3, row 5 -> row 4,
# row 2 -> row 5, row 3 -> row 6, row 6 -> row 7, row 7 -> row 8, row 4 -> row 9.
# Only columns D, L, M, N, O, P, Q, S, T actually change values (the rest
# are identical across all rows already).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44664
$ws.Range("M2").Value = 150
$ws.Range("Q2").Value = "$/caja 18 kilos"
$ws.Range("S2").Value = 1639
$ws.Range("T2").Value = 18

$ws.Range("D3").Value = 44671
$ws.Range("L3").Value = "Segunda"
$ws.Range("N3").Value = 29000
$ws.Range("O3").Value = 30000
$ws.Range("P3").Value = 29500
$ws.Range("S3").Value = 1475

$ws.Range("D4").Value = 44643
$ws.Range("M4").Value = 160
$ws.Range("N4").Value = 28000
$ws.Range("P4").Value = 29000
$ws.Range("S4").Value = 1450

$ws.Range("D5").Value = 44679
$ws.Range("L5").Value = "Segunda"
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 29000
$ws.Range("P5").Value = 29500
$ws.Range("S5").Value = 1475

$ws.Range("D6").Value = 44679
$ws.Range("L6").Value = "Tercera"
$ws.Range("M6").Value = 200
$ws.Range("N6").Value = 24000
$ws.Range("O6").Value = 25000
$ws.Range("P6").Value = 24500
$ws.Range("S6").Value = 1225

$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 160
$ws.Range("N7").Value = 31000
$ws.Range("O7").Value = 32000
$ws.Range("P7").Value = 31500
$ws.Range("S7").Value = 1575

$ws.Range("D8").Value = 44650
$ws.Range("M8").Value = 250

$ws.Range("D9").Value = 44636
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 200
$ws.Range("Q9").Value = "$/caja 20 kilos"
$ws.Range("S9").Value = 1475
$ws.Range("T9").Value = 20
